$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write cell values in the same order the original author typed them so
# that new shared-string entries land at the expected indices (9..16).
$ws.Range("A4").Value = "Training3"
$ws.Range("B4").Value = "-"
$ws.Range("A5").Value = "Training 4"
$ws.Range("C4").Value = "23k"
$ws.Range("D1").Value = "Notes"
$ws.Range("C5").Value = "276k"
$ws.Range("B5").Value = "Espersegueixen, cara a cara, s'ataquen a vegades, escut no utilitzat"
$ws.Range("D5").Value = "Recompenses per distancia mes petites, canviar recompenses de velocitat I escut"

# --- Formatting ---
# D1 "Notes" header: same fill+center style as B1
$ws.Range("B1").Copy()
$ws.Range("D1").PasteSpecial(-4122)  # xlPasteFormats

# C4 / C5 "Steps" values: same right-aligned, no-fill style as C2/C3
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats

# New column D width
$ws.Columns.Item(4).ColumnWidth = 78.1

$excel.CutCopyMode = $false

# --- Update active selection to D6 ---
$ws.Range("D6").Select()
